{"js": "// Replace the quiz answer text in each table cell with the new values,\n// per the recorded diff (24 of the 26 equation cells change; the title\n// paragraph and the \"70\u00f72=35, 0\" cell are left untouched).\nconst replacements = [\n  [\"54\u00f78=6, 6\", \"76\u00f78=9, 4\"],\n  [\"20\u00f74=5, 0\", \"26\u00f79=2, 8\"],\n  [\"32\u00f76=5, 2\", \"64\u00f79=7, 1\"],\n  [\"63\u00f73=21, 0\", \"16\u00f72=8, 0\"],\n  [\"89\u00f77=12, 5\", \"68\u00f75=13, 3\"],\n  [\"87\u00f74=21, 3\", \"74\u00f75=14, 4\"],\n  [\"75\u00f73=25, 0\", \"88\u00f77=12, 4\"],\n  [\"33\u00f73=11, 0\", \"96\u00f74=24, 0\"],\n  [\"68\u00f77=9, 5\", \"69\u00f78=8, 5\"],\n  [\"13\u00f76=2, 1\", \"93\u00f75=18, 3\"],\n  [\"97\u00f73=32, 1\", \"61\u00f73=20, 1\"],\n  [\"70\u00f75=14, 0\", \"86\u00f74=21, 2\"],\n  [\"64\u00f74=16, 0\", \"85\u00f77=12, 1\"],\n  [\"77\u00f74=19, 1\", \"55\u00f75=11, 0\"],\n  [\"91\u00f75=18, 1\", \"45\u00f73=15, 0\"],\n  [\"56\u00f74=14, 0\", \"76\u00f72=38, 0\"],\n  [\"31\u00f72=15, 1\", \"21\u00f77=3, 0\"],\n  [\"41\u00f74=10, 1\", \"51\u00f75=10, 1\"],\n  [\"93\u00f72=46, 1\", \"63\u00f76=10, 3\"],\n  [\"87\u00f72=43, 1\", \"34\u00f76=5, 4\"],\n  [\"89\u00f78=11, 1\", \"34\u00f77=4, 6\"],\n  [\"60\u00f79=6, 6\", \"69\u00f74=17, 1\"],\n  [\"75\u00f79=8, 3\", \"81\u00f79=9, 0\"],\n  [\"42\u00f79=4, 6\", \"63\u00f78=7, 7\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the quiz answer text in each table cell with the new values,\n# per the recorded diff (24 of the 26 equation cells change; the title\n# paragraph and the \"70\u00f72=35, 0\" cell are left untouched).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"54\u00f78=6, 6\", \"76\u00f78=9, 4\"),\n    @(\"20\u00f74=5, 0\", \"26\u00f79=2, 8\"),\n    @(\"32\u00f76=5, 2\", \"64\u00f79=7, 1\"),\n    @(\"63\u00f73=21, 0\", \"16\u00f72=8, 0\"),\n    @(\"89\u00f77=12, 5\", \"68\u00f75=13, 3\"),\n    @(\"87\u00f74=21, 3\", \"74\u00f75=14, 4\"),\n    @(\"75\u00f73=25, 0\", \"88\u00f77=12, 4\"),\n    @(\"33\u00f73=11, 0\", \"96\u00f74=24, 0\"),\n    @(\"68\u00f77=9, 5\", \"69\u00f78=8, 5\"),\n    @(\"13\u00f76=2, 1\", \"93\u00f75=18, 3\"),\n    @(\"97\u00f73=32, 1\", \"61\u00f73=20, 1\"),\n    @(\"70\u00f75=14, 0\", \"86\u00f74=21, 2\"),\n    @(\"64\u00f74=16, 0\", \"85\u00f77=12, 1\"),\n    @(\"77\u00f74=19, 1\", \"55\u00f75=11, 0\"),\n    @(\"91\u00f75=18, 1\", \"45\u00f73=15, 0\"),\n    @(\"56\u00f74=14, 0\", \"76\u00f72=38, 0\"),\n    @(\"31\u00f72=15, 1\", \"21\u00f77=3, 0\"),\n    @(\"41\u00f74=10, 1\", \"51\u00f75=10, 1\"),\n    @(\"93\u00f72=46, 1\", \"63\u00f76=10, 3\"),\n    @(\"87\u00f72=43, 1\", \"34\u00f76=5, 4\"),\n    @(\"89\u00f78=11, 1\", \"34\u00f77=4, 6\"),\n    @(\"60\u00f79=6, 6\", \"69\u00f74=17, 1\"),\n    @(\"75\u00f79=8, 3\", \"81\u00f79=9, 0\"),\n    @(\"42\u00f79=4, 6\", \"63\u00f78=7, 7\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2)\n}\n"}
